$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) figures.
# D-column prices are stored as literal text (to preserve exact
# formatting such as trailing zeros / multi-dot grouping), so we
# force the NumberFormat to Text before writing, then restore the
# cell's style so no stray formatting is left behind.

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.091.98'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +2.46%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.406.48'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +2.95%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '561.03'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +2.31%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '138.21'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +4.79%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.586'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.74%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.404.07'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.93%  '
$ws.Cells.Item(10, 5).Value = '  +2.82%  '
$ws.Cells.Item(11, 5).Value = '  +3.57%  '
$ws.Cells.Item(12, 5).Value = '  +0.05%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.349'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +3.34%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.80'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +8.06%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.834.69'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.91%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.042.72'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +2.40%  '
$ws.Cells.Item(17, 5).Value = '  +3.72%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.398.55'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +2.65%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.03'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +3.49%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '344.19'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +8.81%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.23'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.44%  '
$ws.Cells.Item(22, 5).Value = '  +3.30%  '
$ws.Cells.Item(23, 5).Value = '  +0.36%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '65.07'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.29%  '
$ws.Cells.Item(25, 5).Value = '  +1.26%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.01%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.35'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +5.32%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.51'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +11.09%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.36'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +12.92%  '
$ws.Cells.Item(30, 5).Value = '  +3.64%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0770'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +4.34%  '
$ws.Cells.Item(32, 5).Value = '  +7.03%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '171.70'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.99%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.41'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.17%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.394'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.43%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.53'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +3.25%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.52'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +10.11%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '365.34'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +10.68%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.01%  '
$ws.Cells.Item(40, 5).Value = '  -0.08%  '
$ws.Cells.Item(41, 5).Value = '  +8.19%  '
$ws.Cells.Item(42, 5).Value = '  +2.42%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '143.95'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +2.66%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.67'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +4.87%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.49'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +5.78%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0965'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.93%  '
$ws.Cells.Item(47, 5).Value = '  +4.02%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.583'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.46%  '
$ws.Cells.Item(49, 5).Value = '  +3.40%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.86'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +4.47%  '
$ws.Cells.Item(51, 5).Value = '  -2.83%  '
